# Daily Scores update — 2025-02-01 workbook
# 1) Corrects the existing 2025-02-20 "abs_activity" row (row 78: C/D/F values).
# 2) Appends three new days of scores (2025-02-21, 2025-02-22, 2025-02-23),
#    each with the usual abs_activity / rel_activity / abs_sleep / rel_sleep
#    quadruplet, in rows 82-93.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Fix up existing row 78 (2025-02-20, abs_activity) ------------------
$ws.Cells.Item(78, 3).Value = 9.717843849060959
$ws.Cells.Item(78, 4).Value = 7.92017790036887
$ws.Cells.Item(78, 6).Value = 17.63802174942983

# --- 2) Append the new rows -------------------------------------------------
$newRows = @(
    @{ Row = 82; Date = "2025-02-21"; Typ = "abs_activity"; C = 10;                D = 7.364862235903361;  E = 0; F = 17.36486223590336 }
    @{ Row = 83; Date = "2025-02-21"; Typ = "rel_activity"; C = 5.614872685185185; D = 0;                  E = 0; F = 5.614872685185185 }
    @{ Row = 84; Date = "2025-02-21"; Typ = "abs_sleep";    C = 2.466666666666665; D = 10;                 E = 0; F = 12.46666666666667 }
    @{ Row = 85; Date = "2025-02-21"; Typ = "rel_sleep";    C = 0;                 D = 9.018057880889739;  E = 0; F = 9.018057880889739 }
    @{ Row = 86; Date = "2025-02-22"; Typ = "abs_activity"; C = 9.883803352612128; D = 8.689684777358082;  E = 0; F = 18.57348812997021 }
    @{ Row = 87; Date = "2025-02-22"; Typ = "rel_activity"; C = 0;                 D = 0;                  E = 0; F = 0 }
    @{ Row = 88; Date = "2025-02-22"; Typ = "abs_sleep";    C = 10;                D = 1.466666666666665;  E = 0; F = 11.46666666666667 }
    @{ Row = 89; Date = "2025-02-22"; Typ = "rel_sleep";    C = 10;                D = 0;                  E = 0; F = 10 }
    @{ Row = 90; Date = "2025-02-23"; Typ = "abs_activity"; C = 8.930160447946134; D = 0;                  E = 0; F = 8.930160447946134 }
    @{ Row = 91; Date = "2025-02-23"; Typ = "rel_activity"; C = 0;                 D = 0;                  E = 0; F = 0 }
    @{ Row = 92; Date = "2025-02-23"; Typ = "abs_sleep";    C = 10;                D = 0;                  E = 0; F = 10 }
    @{ Row = 93; Date = "2025-02-23"; Typ = "rel_sleep";    C = 8.648049807727523; D = 0;                  E = 0; F = 8.648049807727523 }
)

foreach ($row in $newRows) {
    $r = $row.Row

    # Format column A as text first so Excel stores the literal date string
    # (e.g. "2025-02-21") instead of auto-converting it to a date serial.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row.Date

    $ws.Cells.Item($r, 2).Value = $row.Typ
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
}
